$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2727.2727
$ws.Range("I51").Value = 2833.3333
$ws.Range("J51").Value = 2687.5
$ws.Range("K51").Value = 2833.3333
$ws.Range("L51").Value = 2687.5
$ws.Range("M51").Value = -2349.3333
$ws.Range("N51").Value = -3655.5

$ws.Range("H108").Value = 30000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 30000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws.Range("H109").Value = 66266.664
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 66266.664
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 66266.664
$ws.Range("N109").Value = -69040.664

$ws.Range("H137").Value = 2571.4614
$ws.Range("I137").Value = 2997
$ws.Range("J137").Value = 2020.7646
$ws.Range("K137").Value = 8991
$ws.Range("L137").Value = 6062.293799999999
$ws.Range("M137").Value = -6441
$ws.Range("N137").Value = -11162.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1682.9286
$ws.Range("I2").Value = 1681.6154
$ws.Range("J2").Value = 1700
$ws.Range("K2").Value = 1681.6154
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = -1568.6154
$ws.Range("N2").Value = -1926

$ws.Range("H32").Value = 18525.535
$ws.Range("I32").Value = 20852.295
$ws.Range("J32").Value = 1573.4286
$ws.Range("K32").Value = 20852.295
$ws.Range("L32").Value = 1573.4286
$ws.Range("M32").Value = -20565.295
$ws.Range("N32").Value = -2147.4286

$ws.Range("H61").Value = 10947.921
$ws.Range("I61").Value = 7779.1035
$ws.Range("J61").Value = 21158.555
$ws.Range("K61").Value = 7779.1035
$ws.Range("L61").Value = 21158.555
$ws.Range("M61").Value = -7567.1035
$ws.Range("N61").Value = -21582.555

$ws.Range("H116").Value = 1682.9286
$ws.Range("I116").Value = 1681.6154
$ws.Range("J116").Value = 1700
$ws.Range("K116").Value = 1681.6154
$ws.Range("L116").Value = 1700
$ws.Range("M116").Value = 612.3846000000001
$ws.Range("N116").Value = -6288

$ws.Range("H132").Value = 2730.2273
$ws.Range("I132").Value = 1931.3572
$ws.Range("J132").Value = 4128.25
$ws.Range("K132").Value = 5794.071599999999
$ws.Range("L132").Value = 12384.75
$ws.Range("M132").Value = -3264.071599999999
$ws.Range("N132").Value = -17444.75

$ws.Range("H136").Value = 10947.921
$ws.Range("I136").Value = 7779.1035
$ws.Range("J136").Value = 21158.555
$ws.Range("K136").Value = 23337.3105
$ws.Range("L136").Value = 63475.665
$ws.Range("M136").Value = -20787.3105
$ws.Range("N136").Value = -68575.66500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1682.9286
$ws.Range("I3").Value = 1681.6154
$ws.Range("J3").Value = 1700
$ws.Range("K3").Value = 1681.6154
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = -1567.6154
$ws.Range("N3").Value = -1928

$ws.Range("H18").Value = 19000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 19000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 19000
$ws.Range("N18").Value = -20058

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5418.1113
$ws.Range("I31").Value = 4566.1934
$ws.Range("J31").Value = 10700
$ws.Range("K31").Value = 4566.1934
$ws.Range("L31").Value = 10700
$ws.Range("M31").Value = -4271.1934
$ws.Range("N31").Value = -11290

$ws.Range("H34").Value = 5418.1113
$ws.Range("I34").Value = 4566.1934
$ws.Range("J34").Value = 10700
$ws.Range("K34").Value = 4566.1934
$ws.Range("L34").Value = 10700
$ws.Range("M34").Value = -4364.1934
$ws.Range("N34").Value = -11104

$ws.Range("H58").Value = 2527836
$ws.Range("I58").Value = 4786837.5
$ws.Range("J58").Value = 3069.647
$ws.Range("K58").Value = 4786837.5
$ws.Range("L58").Value = 3069.647
$ws.Range("M58").Value = -4786634.5
$ws.Range("N58").Value = -3475.647

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H132").Value = 3278.9473
$ws.Range("I132").Value = 2553.5334
$ws.Range("J132").Value = 5999.25
$ws.Range("K132").Value = 7660.600199999999
$ws.Range("L132").Value = 17997.75
$ws.Range("M132").Value = -5130.600199999999
$ws.Range("N132").Value = -23057.75

$ws.Range("H136").Value = 2527836
$ws.Range("I136").Value = 4786837.5
$ws.Range("J136").Value = 3069.647
$ws.Range("K136").Value = 14360512.5
$ws.Range("L136").Value = 9208.940999999999
$ws.Range("M136").Value = -14357962.5
$ws.Range("N136").Value = -14308.941

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8338771
$ws.Range("I5").Value = 358.75
$ws.Range("J5").Value = 41692420
$ws.Range("K5").Value = 1076.25
$ws.Range("L5").Value = 125077260
$ws.Range("M5").Value = -964.25
$ws.Range("N5").Value = -125077484

$ws.Range("H10").Value = 849.5
$ws.Range("I10").Value = 849.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 2548.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -2409.5

$ws.Range("H113").Value = 689.7857
$ws.Range("I113").Value = 691.7778
$ws.Range("J113").Value = 683.0625
$ws.Range("K113").Value = 2075.3334
$ws.Range("L113").Value = 2049.1875
$ws.Range("M113").Value = 94.66660000000002
$ws.Range("N113").Value = -6389.1875

$ws.Range("H122").Value = 1024.7059
$ws.Range("I122").Value = 369.33334
$ws.Range("J122").Value = 1382.1818
$ws.Range("K122").Value = 3324.00006
$ws.Range("L122").Value = 12439.6362
$ws.Range("M122").Value = -874.0000600000003
$ws.Range("N122").Value = -17339.6362

$ws.Range("H132").Value = 1350.3077
$ws.Range("I132").Value = 919.25
$ws.Range("J132").Value = 2040
$ws.Range("K132").Value = 8273.25
$ws.Range("L132").Value = 18360
$ws.Range("M132").Value = -5743.25
$ws.Range("N132").Value = -23420

$ws.Range("H135").Value = 8338771
$ws.Range("I135").Value = 358.75
$ws.Range("J135").Value = 41692420
$ws.Range("K135").Value = 3228.75
$ws.Range("L135").Value = 375231780
$ws.Range("M135").Value = -693.75
$ws.Range("N135").Value = -375236850

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10928.571
$ws.Range("I5").Value = 6000
$ws.Range("J5").Value = 17500
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 17500
$ws.Range("M5").Value = -5888
$ws.Range("N5").Value = -17724

$ws.Range("H80").Value = 5798.591
$ws.Range("I80").Value = 7700.9
$ws.Range("J80").Value = 4213.3335
$ws.Range("K80").Value = 7700.9
$ws.Range("L80").Value = 4213.3335
$ws.Range("M80").Value = -6702.9
$ws.Range("N80").Value = -6209.3335

$ws.Range("H83").Value = 5798.591
$ws.Range("I83").Value = 7700.9
$ws.Range("J83").Value = 4213.3335
$ws.Range("K83").Value = 38504.5
$ws.Range("L83").Value = 21066.6675
$ws.Range("M83").Value = -33512.5
$ws.Range("N83").Value = -31050.6675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5676.1665
$ws.Range("I132").Value = 6202.2
$ws.Range("J132").Value = 5300.4287
$ws.Range("K132").Value = 18606.6
$ws.Range("L132").Value = 15901.2861
$ws.Range("M132").Value = -16076.6
$ws.Range("N132").Value = -20961.2861

$ws.Range("H134").Value = 43209.668
$ws.Range("I134").Value = 20000
$ws.Range("J134").Value = 54814.5
$ws.Range("K134").Value = 20000
$ws.Range("L134").Value = 54814.5
$ws.Range("M134").Value = -14930
$ws.Range("N134").Value = -64954.5
